$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.050.91"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.813.62"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("D4").Value = "'0.9976"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'310.76"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "'0.9982"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'0.4417"
$ws.Range("E7").Value = "  +4.49%  "
$ws.Range("D8").Value = "'0.3719"
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("D9").Value = "'0.07441"
$ws.Range("E9").Value = "  +2.98%  "
$ws.Range("D10").Value = "'0.8665"
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("D11").Value = "'20.77"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "1.798.20"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").Value = "'6.666"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "'93.34"
$ws.Range("E14").Value = "  +3.95%  "
$ws.Range("D15").Value = "'0.07073"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "'5.294"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "'0.9984"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "'0.000008714"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").Value = "'0.9988"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "'14.88"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "27.041.01"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").Value = "'5.184"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").Value = "'10.86"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").Value = "2.023.49"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").Value = "'1.976"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").Value = "'151.14"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").Value = "'2.220"
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("D28").Value = "'18.40"
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("D29").Value = "'5.222"
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("D30").Value = "'117.52"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").Value = "'0.08788"
$ws.Range("E31").Value = "  +0.62%  "
$ws.Range("D32").Value = "'0.7503"
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").Value = "'1.171"
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("D34").Value = "'4.504"
$ws.Range("E34").Value = "  +1.82%  "
$ws.Range("D35").Value = "'2.884"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "'0.9970"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").Value = "'1.095"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").Value = "'0.01974"
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("D39").Value = "'0.05231"
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").Value = "'0.5267"
$ws.Range("E40").Value = "  +3.91%  "
$ws.Range("D41").Value = "'7.132"
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("D42").Value = "'2.824"
$ws.Range("E42").Value = "  -1.84%  "
$ws.Range("D43").Value = "'0.1695"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").Value = "'2.178"
$ws.Range("E44").Value = "  +11.64%  "
$ws.Range("D45").Value = "'8.567"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'0.4999"
$ws.Range("E46").Value = "  +5.73%  "
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("D48").Value = "'104.61"
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("D49").Value = "'1.680"
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").Value = "'0.9980"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").Value = "'0.06340"
$ws.Range("E51").Value = "  +0.21%  "
